$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: "Contact types" ---
$ws.Range("A5").Value = "Contact types"
$ws.Range("B5").Value = "Contacts should be able to be of zero to many types. Currently they were only able to be one type."
$ws.Range("C5").Value = "OPEN"
$ws.Rows.Item(5).RowHeight = 30

# --- Row 6: "Drag contacts between addresses" ---
$ws.Range("A6").Value = "Drag contacts between addresses"
$ws.Range("B6").Value = "In the Addresses and Contacts screen it should be possible to drag a contact from one address into another."
$ws.Range("C6").Value = "OPEN"
$ws.Rows.Item(6).RowHeight = 30

# --- Row 7: "Use new AutoCompleteBox" (description entered before title, matching original authoring order) ---
$ws.Range("B7").Value = "Consider replacing comboboxes with text searches with the new AutoCompleteBox control that is included in the newest WPF Toolkit release. Docs at http://msdn.microsoft.com/en-us/library/system.windows.controls.autocompletebox(VS.95).aspx"
$ws.Range("A7").Value = "Use new AutoCompleteBox"
$ws.Range("C7").Value = "OPEN"
$ws.Rows.Item(7).RowHeight = 45

# Move the active selection to A7, matching the saved view state after entry
$ws.Range("A7").Select() | Out-Null
